$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.968.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.786.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +22.71%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.79"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +7.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.67"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.782.29"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +22.65%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.543"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.13%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +10.34%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.500"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.83"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +14.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000258"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.417.20"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +22.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.789.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +22.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.125.63"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +6.19%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "522.93"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +7.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.73"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.96%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +23.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.747"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +9.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.48"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.29%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +11.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.57"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.04"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +9.09%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000122"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +29.67%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +10.57%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +13.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.03"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.24"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +15.52%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +12.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.15"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +10.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.23"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +11.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.340"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.91%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.39"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.161.09"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +13.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "430.82"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +17.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.84"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +7.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.04"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.62%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.31%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.82"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.68"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.48"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.80%  "
